$wb = $excel.ActiveWorkbook

# --- Sheet 1 ---
$ws = $wb.Worksheets.Item(1)
$ws.Name = 'summ38851398'
$ws.Rows.Item(4).Delete()
$ws.Cells.Item(2,1).Value = 'Intercept'
$ws.Cells.Item(2,2).Value = [double]"8559.4792142925"
$ws.Cells.Item(2,3).Value = [double]"2.104703365986733e-10"
$ws.Cells.Item(3,1).Value = 'Education[T.Secondary]'
$ws.Cells.Item(3,2).Value = [double]"85.56644859213452"
$ws.Cells.Item(3,3).Value = [double]"0.9100168092499552"
$ws.Cells.Item(4,1).Value = 'Education[T.University]'
$ws.Cells.Item(4,2).Value = [double]"1038.474093382618"
$ws.Cells.Item(4,3).Value = [double]"0.1684449262855488"
$ws.Cells.Item(5,1).Value = 'Education[T.Unknown/Other]'
$ws.Cells.Item(5,2).Value = [double]"-1919.830912901279"
$ws.Cells.Item(5,3).Value = [double]"0.1348018027623009"
$ws.Cells.Item(6,1).Value = 'HHSize'
$ws.Cells.Item(6,2).Value = [double]"185.740277792478"
$ws.Cells.Item(6,3).Value = [double]"0.01732532415407506"
$ws.Cells.Item(7,1).Value = 'Sex'
$ws.Cells.Item(7,2).Value = [double]"-1508.61272696468"
$ws.Cells.Item(7,3).Value = [double]"2.589641352173449e-15"
$ws.Cells.Item(8,1).Value = 'Age'
$ws.Cells.Item(8,2).Value = [double]"-36.32249078929323"
$ws.Cells.Item(8,3).Value = [double]"1.467628002996332e-05"
$ws.Cells.Item(9,1).Value = 'DistSubcenter_res'
$ws.Cells.Item(9,2).Value = [double]"431.1455622058448"
$ws.Cells.Item(9,3).Value = [double]"1.175101843345304e-16"
$ws.Cells.Item(10,1).Value = 'DistCenter_res'
$ws.Cells.Item(10,2).Value = [double]"225.1144994421358"
$ws.Cells.Item(10,3).Value = [double]"1.073217817351589e-09"
$ws.Cells.Item(11,1).Value = 'UrbPopDensity_res'
$ws.Cells.Item(11,2).Value = [double]"0.08638309601867286"
$ws.Cells.Item(11,3).Value = [double]"0.228526137107308"
$ws.Cells.Item(12,1).Value = 'UrbBuildDensity_res'
$ws.Cells.Item(12,2).Value = [double]"-0.0001370410827594626"
$ws.Cells.Item(12,3).Value = [double]"0.02578860560074952"
$ws.Cells.Item(13,1).Value = 'IntersecDensity_res'
$ws.Cells.Item(13,2).Value = [double]"-0.1109709129630811"
$ws.Cells.Item(13,3).Value = [double]"0.9873237093736872"
$ws.Cells.Item(14,1).Value = 'street_length_res'
$ws.Cells.Item(14,2).Value = [double]"24.59662067461386"
$ws.Cells.Item(14,3).Value = [double]"5.202620712388441e-06"
$ws.Cells.Item(15,1).Value = 'LU_Comm_res'
$ws.Cells.Item(15,2).Value = [double]"-2248.751920221544"
$ws.Cells.Item(15,3).Value = [double]"0.04462858924189807"
$ws.Cells.Item(16,1).Value = 'LU_UrbFab_res'
$ws.Cells.Item(16,2).Value = [double]"-2987.105194089568"
$ws.Cells.Item(16,3).Value = [double]"2.074905062507086e-05"
$ws.Cells.Item(17,1).Value = 'bike_lane_share_res'
$ws.Cells.Item(17,2).Value = [double]"-2326.967088740204"
$ws.Cells.Item(17,3).Value = [double]"0.1665771482867518"

# --- Sheet 2 ---
$ws = $wb.Worksheets.Item(2)
$ws.Name = 'summ39217429'
$ws.Rows.Item(4).Delete()
$ws.Cells.Item(2,1).Value = 'Intercept'
$ws.Cells.Item(2,2).Value = [double]"8309.14308277113"
$ws.Cells.Item(2,3).Value = [double]"1.62237491557614e-09"
$ws.Cells.Item(3,1).Value = 'Education[T.Secondary]'
$ws.Cells.Item(3,2).Value = [double]"1094.805624574649"
$ws.Cells.Item(3,3).Value = [double]"0.1590588593882504"
$ws.Cells.Item(4,1).Value = 'Education[T.University]'
$ws.Cells.Item(4,2).Value = [double]"1959.281149070185"
$ws.Cells.Item(4,3).Value = [double]"0.01147282100397223"
$ws.Cells.Item(5,1).Value = 'Education[T.Unknown/Other]'
$ws.Cells.Item(5,2).Value = [double]"-677.4080344482131"
$ws.Cells.Item(5,3).Value = [double]"0.5667766273750298"
$ws.Cells.Item(6,1).Value = 'HHSize'
$ws.Cells.Item(6,2).Value = [double]"234.5114139607753"
$ws.Cells.Item(6,3).Value = [double]"0.002841903028612641"
$ws.Cells.Item(7,1).Value = 'Sex'
$ws.Cells.Item(7,2).Value = [double]"-1783.81269440894"
$ws.Cells.Item(7,3).Value = [double]"7.026446231848277e-21"
$ws.Cells.Item(8,1).Value = 'Age'
$ws.Cells.Item(8,2).Value = [double]"-39.10397731969555"
$ws.Cells.Item(8,3).Value = [double]"3.875183611787796e-06"
$ws.Cells.Item(9,1).Value = 'DistSubcenter_res'
$ws.Cells.Item(9,2).Value = [double]"463.6332288026471"
$ws.Cells.Item(9,3).Value = [double]"1.484227295772252e-18"
$ws.Cells.Item(10,1).Value = 'DistCenter_res'
$ws.Cells.Item(10,2).Value = [double]"213.6606185143538"
$ws.Cells.Item(10,3).Value = [double]"8.277000203851175e-09"
$ws.Cells.Item(11,1).Value = 'UrbPopDensity_res'
$ws.Cells.Item(11,2).Value = [double]"0.1002064085749313"
$ws.Cells.Item(11,3).Value = [double]"0.155649489870486"
$ws.Cells.Item(12,1).Value = 'UrbBuildDensity_res'
$ws.Cells.Item(12,2).Value = [double]"-0.0001844164623483085"
$ws.Cells.Item(12,3).Value = [double]"0.003156211319440929"
$ws.Cells.Item(13,1).Value = 'IntersecDensity_res'
$ws.Cells.Item(13,2).Value = [double]"0.01193919792281317"
$ws.Cells.Item(13,3).Value = [double]"0.99863052723613"
$ws.Cells.Item(14,1).Value = 'street_length_res'
$ws.Cells.Item(14,2).Value = [double]"25.24264602736938"
$ws.Cells.Item(14,3).Value = [double]"6.544573652570229e-06"
$ws.Cells.Item(15,1).Value = 'LU_Comm_res'
$ws.Cells.Item(15,2).Value = [double]"-2984.867281723178"
$ws.Cells.Item(15,3).Value = [double]"0.007553186463040676"
$ws.Cells.Item(16,1).Value = 'LU_UrbFab_res'
$ws.Cells.Item(16,2).Value = [double]"-3581.448139510693"
$ws.Cells.Item(16,3).Value = [double]"4.250138601556027e-07"
$ws.Cells.Item(17,1).Value = 'bike_lane_share_res'
$ws.Cells.Item(17,2).Value = [double]"-1145.129625431258"
$ws.Cells.Item(17,3).Value = [double]"0.4914736723823991"

# --- Sheet 3 ---
$ws = $wb.Worksheets.Item(3)
$ws.Name = 'summ39628127'
$ws.Rows.Item(4).Delete()
$ws.Cells.Item(2,1).Value = 'Intercept'
$ws.Cells.Item(2,2).Value = [double]"8931.500701425597"
$ws.Cells.Item(2,3).Value = [double]"5.22659381760032e-11"
$ws.Cells.Item(3,1).Value = 'Education[T.Secondary]'
$ws.Cells.Item(3,2).Value = [double]"496.4167703936287"
$ws.Cells.Item(3,3).Value = [double]"0.5265697602032553"
$ws.Cells.Item(4,1).Value = 'Education[T.University]'
$ws.Cells.Item(4,2).Value = [double]"1352.839206589945"
$ws.Cells.Item(4,3).Value = [double]"0.08301988411563067"
$ws.Cells.Item(5,1).Value = 'Education[T.Unknown/Other]'
$ws.Cells.Item(5,2).Value = [double]"-1449.380677945643"
$ws.Cells.Item(5,3).Value = [double]"0.216055056443741"
$ws.Cells.Item(6,1).Value = 'HHSize'
$ws.Cells.Item(6,2).Value = [double]"205.5087619567472"
$ws.Cells.Item(6,3).Value = [double]"0.008716459025500847"
$ws.Cells.Item(7,1).Value = 'Sex'
$ws.Cells.Item(7,2).Value = [double]"-1644.773959870839"
$ws.Cells.Item(7,3).Value = [double]"4.891141269586353e-18"
$ws.Cells.Item(8,1).Value = 'Age'
$ws.Cells.Item(8,2).Value = [double]"-45.60125572321082"
$ws.Cells.Item(8,3).Value = [double]"5.788019384202165e-08"
$ws.Cells.Item(9,1).Value = 'DistSubcenter_res'
$ws.Cells.Item(9,2).Value = [double]"450.0201426670487"
$ws.Cells.Item(9,3).Value = [double]"5.539630246177607e-18"
$ws.Cells.Item(10,1).Value = 'DistCenter_res'
$ws.Cells.Item(10,2).Value = [double]"258.4231545297431"
$ws.Cells.Item(10,3).Value = [double]"3.072641871714979e-12"
$ws.Cells.Item(11,1).Value = 'UrbPopDensity_res'
$ws.Cells.Item(11,2).Value = [double]"0.1505866782085918"
$ws.Cells.Item(11,3).Value = [double]"0.03430709222391803"
$ws.Cells.Item(12,1).Value = 'UrbBuildDensity_res'
$ws.Cells.Item(12,2).Value = [double]"-0.0001472130903816358"
$ws.Cells.Item(12,3).Value = [double]"0.01354641399318221"
$ws.Cells.Item(13,1).Value = 'IntersecDensity_res'
$ws.Cells.Item(13,2).Value = [double]"3.336972871115765"
$ws.Cells.Item(13,3).Value = [double]"0.631143160514041"
$ws.Cells.Item(14,1).Value = 'street_length_res'
$ws.Cells.Item(14,2).Value = [double]"22.33888436758924"
$ws.Cells.Item(14,3).Value = [double]"3.853072953585204e-05"
$ws.Cells.Item(15,1).Value = 'LU_Comm_res'
$ws.Cells.Item(15,2).Value = [double]"-2849.518211552358"
$ws.Cells.Item(15,3).Value = [double]"0.009798753278414639"
$ws.Cells.Item(16,1).Value = 'LU_UrbFab_res'
$ws.Cells.Item(16,2).Value = [double]"-4021.181109881573"
$ws.Cells.Item(16,3).Value = [double]"9.721755231894737e-09"
$ws.Cells.Item(17,1).Value = 'bike_lane_share_res'
$ws.Cells.Item(17,2).Value = [double]"-2798.776755381859"
$ws.Cells.Item(17,3).Value = [double]"0.09858169099027267"

# --- Sheet 4 ---
$ws = $wb.Worksheets.Item(4)
$ws.Name = 'summ40061656'
$ws.Rows.Item(4).Delete()
$ws.Cells.Item(2,1).Value = 'Intercept'
$ws.Cells.Item(2,2).Value = [double]"8806.796907265918"
$ws.Cells.Item(2,3).Value = [double]"1.430810888449753e-10"
$ws.Cells.Item(3,1).Value = 'Education[T.Secondary]'
$ws.Cells.Item(3,2).Value = [double]"1346.803444201805"
$ws.Cells.Item(3,3).Value = [double]"0.0944648399965085"
$ws.Cells.Item(4,1).Value = 'Education[T.University]'
$ws.Cells.Item(4,2).Value = [double]"2347.599304397963"
$ws.Cells.Item(4,3).Value = [double]"0.003452822730656717"
$ws.Cells.Item(5,1).Value = 'Education[T.Unknown/Other]'
$ws.Cells.Item(5,2).Value = [double]"-511.7770640520671"
$ws.Cells.Item(5,3).Value = [double]"0.6699434250648812"
$ws.Cells.Item(6,1).Value = 'HHSize'
$ws.Cells.Item(6,2).Value = [double]"179.049273571018"
$ws.Cells.Item(6,3).Value = [double]"0.02122158937610778"
$ws.Cells.Item(7,1).Value = 'Sex'
$ws.Cells.Item(7,2).Value = [double]"-1509.25977131059"
$ws.Cells.Item(7,3).Value = [double]"1.179108793230331e-15"
$ws.Cells.Item(8,1).Value = 'Age'
$ws.Cells.Item(8,2).Value = [double]"-37.287673230363"
$ws.Cells.Item(8,3).Value = [double]"7.831235737443233e-06"
$ws.Cells.Item(9,1).Value = 'DistSubcenter_res'
$ws.Cells.Item(9,2).Value = [double]"363.9348707257268"
$ws.Cells.Item(9,3).Value = [double]"2.00371285510003e-12"
$ws.Cells.Item(10,1).Value = 'DistCenter_res'
$ws.Cells.Item(10,2).Value = [double]"230.5992288745629"
$ws.Cells.Item(10,3).Value = [double]"3.509957197825063e-10"
$ws.Cells.Item(11,1).Value = 'UrbPopDensity_res'
$ws.Cells.Item(11,2).Value = [double]"0.1024409138053081"
$ws.Cells.Item(11,3).Value = [double]"0.1472046090455995"
$ws.Cells.Item(12,1).Value = 'UrbBuildDensity_res'
$ws.Cells.Item(12,2).Value = [double]"-0.0001498544457121694"
$ws.Cells.Item(12,3).Value = [double]"0.01103951224305827"
$ws.Cells.Item(13,1).Value = 'IntersecDensity_res'
$ws.Cells.Item(13,2).Value = [double]"-4.432226811669192"
$ws.Cells.Item(13,3).Value = [double]"0.5206139873059128"
$ws.Cells.Item(14,1).Value = 'street_length_res'
$ws.Cells.Item(14,2).Value = [double]"17.63778487602703"
$ws.Cells.Item(14,3).Value = [double]"0.0012040025531448"
$ws.Cells.Item(15,1).Value = 'LU_Comm_res'
$ws.Cells.Item(15,2).Value = [double]"-3153.444102883986"
$ws.Cells.Item(15,3).Value = [double]"0.004127478805718771"
$ws.Cells.Item(16,1).Value = 'LU_UrbFab_res'
$ws.Cells.Item(16,2).Value = [double]"-3546.794451671588"
$ws.Cells.Item(16,3).Value = [double]"3.810784273923718e-07"
$ws.Cells.Item(17,1).Value = 'bike_lane_share_res'
$ws.Cells.Item(17,2).Value = [double]"-1996.381649249989"
$ws.Cells.Item(17,3).Value = [double]"0.2323072644887631"

# --- Sheet 5 ---
$ws = $wb.Worksheets.Item(5)
$ws.Name = 'summ40451338'
$ws.Rows.Item(4).Delete()
$ws.Cells.Item(2,1).Value = 'Intercept'
$ws.Cells.Item(2,2).Value = [double]"8485.523670121347"
$ws.Cells.Item(2,3).Value = [double]"7.791850940977217e-10"
$ws.Cells.Item(3,1).Value = 'Education[T.Secondary]'
$ws.Cells.Item(3,2).Value = [double]"617.3672286805429"
$ws.Cells.Item(3,3).Value = [double]"0.427744292432982"
$ws.Cells.Item(4,1).Value = 'Education[T.University]'
$ws.Cells.Item(4,2).Value = [double]"1591.864915374179"
$ws.Cells.Item(4,3).Value = [double]"0.04012669570886505"
$ws.Cells.Item(5,1).Value = 'Education[T.Unknown/Other]'
$ws.Cells.Item(5,2).Value = [double]"74.88671132634954"
$ws.Cells.Item(5,3).Value = [double]"0.9507555515827182"
$ws.Cells.Item(6,1).Value = 'HHSize'
$ws.Cells.Item(6,2).Value = [double]"177.162513169807"
$ws.Cells.Item(6,3).Value = [double]"0.02421344669463003"
$ws.Cells.Item(7,1).Value = 'Sex'
$ws.Cells.Item(7,2).Value = [double]"-1454.744086774182"
$ws.Cells.Item(7,3).Value = [double]"1.815010312762915e-14"
$ws.Cells.Item(8,1).Value = 'Age'
$ws.Cells.Item(8,2).Value = [double]"-47.05751802141307"
$ws.Cells.Item(8,3).Value = [double]"2.267537359155867e-08"
$ws.Cells.Item(9,1).Value = 'DistSubcenter_res'
$ws.Cells.Item(9,2).Value = [double]"432.4057822346945"
$ws.Cells.Item(9,3).Value = [double]"1.444159158537377e-16"
$ws.Cells.Item(10,1).Value = 'DistCenter_res'
$ws.Cells.Item(10,2).Value = [double]"232.176776797333"
$ws.Cells.Item(10,3).Value = [double]"2.8725043465328e-10"
$ws.Cells.Item(11,1).Value = 'UrbPopDensity_res'
$ws.Cells.Item(11,2).Value = [double]"0.06719086712067363"
$ws.Cells.Item(11,3).Value = [double]"0.3442897403884454"
$ws.Cells.Item(12,1).Value = 'UrbBuildDensity_res'
$ws.Cells.Item(12,2).Value = [double]"-0.0001293929374899603"
$ws.Cells.Item(12,3).Value = [double]"0.03062848492933417"
$ws.Cells.Item(13,1).Value = 'IntersecDensity_res'
$ws.Cells.Item(13,2).Value = [double]"1.345430508210871"
$ws.Cells.Item(13,3).Value = [double]"0.8462859845964477"
$ws.Cells.Item(14,1).Value = 'street_length_res'
$ws.Cells.Item(14,2).Value = [double]"25.0159823879047"
$ws.Cells.Item(14,3).Value = [double]"1.15369204280982e-05"
$ws.Cells.Item(15,1).Value = 'LU_Comm_res'
$ws.Cells.Item(15,2).Value = [double]"-2875.0679391975"
$ws.Cells.Item(15,3).Value = [double]"0.009352280930909609"
$ws.Cells.Item(16,1).Value = 'LU_UrbFab_res'
$ws.Cells.Item(16,2).Value = [double]"-3371.828457872727"
$ws.Cells.Item(16,3).Value = [double]"1.354370557360546e-06"
$ws.Cells.Item(17,1).Value = 'bike_lane_share_res'
$ws.Cells.Item(17,2).Value = [double]"-1600.97834558835"
$ws.Cells.Item(17,3).Value = [double]"0.3342061836649873"

# --- Sheet 6 ---
$ws = $wb.Worksheets.Item(6)
$ws.Name = 'summ40855390'
$ws.Rows.Item(4).Delete()
$ws.Cells.Item(2,1).Value = 'Intercept'
$ws.Cells.Item(2,2).Value = [double]"9115.353010299475"
$ws.Cells.Item(2,3).Value = [double]"3.889364591881043e-11"
$ws.Cells.Item(3,1).Value = 'Education[T.Secondary]'
$ws.Cells.Item(3,2).Value = [double]"91.2182943090636"
$ws.Cells.Item(3,3).Value = [double]"0.9095917270772188"
$ws.Cells.Item(4,1).Value = 'Education[T.University]'
$ws.Cells.Item(4,2).Value = [double]"983.6506777193122"
$ws.Cells.Item(4,3).Value = [double]"0.2190306448894407"
$ws.Cells.Item(5,1).Value = 'Education[T.Unknown/Other]'
$ws.Cells.Item(5,2).Value = [double]"-739.814613105344"
$ws.Cells.Item(5,3).Value = [double]"0.5450936830303266"
$ws.Cells.Item(6,1).Value = 'HHSize'
$ws.Cells.Item(6,2).Value = [double]"210.8438177651651"
$ws.Cells.Item(6,3).Value = [double]"0.006979688293958795"
$ws.Cells.Item(7,1).Value = 'Sex'
$ws.Cells.Item(7,2).Value = [double]"-1549.54118773107"
$ws.Cells.Item(7,3).Value = [double]"3.078892094266292e-16"
$ws.Cells.Item(8,1).Value = 'Age'
$ws.Cells.Item(8,2).Value = [double]"-42.22625253362901"
$ws.Cells.Item(8,3).Value = [double]"3.501638114842269e-07"
$ws.Cells.Item(9,1).Value = 'DistSubcenter_res'
$ws.Cells.Item(9,2).Value = [double]"408.9650882066287"
$ws.Cells.Item(9,3).Value = [double]"2.330514139274389e-15"
$ws.Cells.Item(10,1).Value = 'DistCenter_res'
$ws.Cells.Item(10,2).Value = [double]"204.433681199621"
$ws.Cells.Item(10,3).Value = [double]"2.506555849169713e-08"
$ws.Cells.Item(11,1).Value = 'UrbPopDensity_res'
$ws.Cells.Item(11,2).Value = [double]"0.08198519896160054"
$ws.Cells.Item(11,3).Value = [double]"0.2500015776660079"
$ws.Cells.Item(12,1).Value = 'UrbBuildDensity_res'
$ws.Cells.Item(12,2).Value = [double]"-0.0001237558061867247"
$ws.Cells.Item(12,3).Value = [double]"0.04198803009168137"
$ws.Cells.Item(13,1).Value = 'IntersecDensity_res'
$ws.Cells.Item(13,2).Value = [double]"-1.415908191466152"
$ws.Cells.Item(13,3).Value = [double]"0.8378884075182875"
$ws.Cells.Item(14,1).Value = 'street_length_res'
$ws.Cells.Item(14,2).Value = [double]"25.4680967221373"
$ws.Cells.Item(14,3).Value = [double]"2.700056813535874e-06"
$ws.Cells.Item(15,1).Value = 'LU_Comm_res'
$ws.Cells.Item(15,2).Value = [double]"-2685.259660994704"
$ws.Cells.Item(15,3).Value = [double]"0.01636818420562554"
$ws.Cells.Item(16,1).Value = 'LU_UrbFab_res'
$ws.Cells.Item(16,2).Value = [double]"-2983.772227029259"
$ws.Cells.Item(16,3).Value = [double]"2.35000003394536e-05"
$ws.Cells.Item(17,1).Value = 'bike_lane_share_res'
$ws.Cells.Item(17,2).Value = [double]"-2662.151621846941"
$ws.Cells.Item(17,3).Value = [double]"0.1107795939444428"

# --- Sheet 7 ---
$ws = $wb.Worksheets.Item(7)
$ws.Name = 'summ41261600'
$ws.Rows.Item(4).Delete()
$ws.Cells.Item(2,1).Value = 'Intercept'
$ws.Cells.Item(2,2).Value = [double]"11456.62263092298"
$ws.Cells.Item(2,3).Value = [double]"4.890490271394243e-17"
$ws.Cells.Item(3,1).Value = 'Education[T.Secondary]'
$ws.Cells.Item(3,2).Value = [double]"-249.8590082214145"
$ws.Cells.Item(3,3).Value = [double]"0.7477685625579482"
$ws.Cells.Item(4,1).Value = 'Education[T.University]'
$ws.Cells.Item(4,2).Value = [double]"610.604319336249"
$ws.Cells.Item(4,3).Value = [double]"0.4295846642259475"
$ws.Cells.Item(5,1).Value = 'Education[T.Unknown/Other]'
$ws.Cells.Item(5,2).Value = [double]"-1798.515916640818"
$ws.Cells.Item(5,3).Value = [double]"0.1327547524589278"
$ws.Cells.Item(6,1).Value = 'HHSize'
$ws.Cells.Item(6,2).Value = [double]"131.2032564278439"
$ws.Cells.Item(6,3).Value = [double]"0.09315952419990839"
$ws.Cells.Item(7,1).Value = 'Sex'
$ws.Cells.Item(7,2).Value = [double]"-1574.625353065934"
$ws.Cells.Item(7,3).Value = [double]"8.555261857981018e-17"
$ws.Cells.Item(8,1).Value = 'Age'
$ws.Cells.Item(8,2).Value = [double]"-46.92193167789622"
$ws.Cells.Item(8,3).Value = [double]"2.57319721577168e-08"
$ws.Cells.Item(9,1).Value = 'DistSubcenter_res'
$ws.Cells.Item(9,2).Value = [double]"389.4527025677106"
$ws.Cells.Item(9,3).Value = [double]"7.468914910069069e-14"
$ws.Cells.Item(10,1).Value = 'DistCenter_res'
$ws.Cells.Item(10,2).Value = [double]"211.8916694414366"
$ws.Cells.Item(10,3).Value = [double]"1.044486680971026e-08"
$ws.Cells.Item(11,1).Value = 'UrbPopDensity_res'
$ws.Cells.Item(11,2).Value = [double]"0.03617853288152988"
$ws.Cells.Item(11,3).Value = [double]"0.6112909147623342"
$ws.Cells.Item(12,1).Value = 'UrbBuildDensity_res'
$ws.Cells.Item(12,2).Value = [double]"-7.526117392886085e-05"
$ws.Cells.Item(12,3).Value = [double]"0.2346437821305851"
$ws.Cells.Item(13,1).Value = 'IntersecDensity_res'
$ws.Cells.Item(13,2).Value = [double]"-3.574824635130463"
$ws.Cells.Item(13,3).Value = [double]"0.6084127295232047"
$ws.Cells.Item(14,1).Value = 'street_length_res'
$ws.Cells.Item(14,2).Value = [double]"17.25488277075075"
$ws.Cells.Item(14,3).Value = [double]"0.001543090236456346"
$ws.Cells.Item(15,1).Value = 'LU_Comm_res'
$ws.Cells.Item(15,2).Value = [double]"-4525.000438661249"
$ws.Cells.Item(15,3).Value = [double]"4.589903444257186e-05"
$ws.Cells.Item(16,1).Value = 'LU_UrbFab_res'
$ws.Cells.Item(16,2).Value = [double]"-3754.565789107786"
$ws.Cells.Item(16,3).Value = [double]"6.63038171530003e-08"
$ws.Cells.Item(17,1).Value = 'bike_lane_share_res'
$ws.Cells.Item(17,2).Value = [double]"-2491.455904664843"
$ws.Cells.Item(17,3).Value = [double]"0.1335628416931118"

# --- Sheet 8 ---
$ws = $wb.Worksheets.Item(8)
$ws.Name = 'summ41660659'
$ws.Rows.Item(4).Delete()
$ws.Cells.Item(2,1).Value = 'Intercept'
$ws.Cells.Item(2,2).Value = [double]"8264.237077695803"
$ws.Cells.Item(2,3).Value = [double]"1.93698835288573e-09"
$ws.Cells.Item(3,1).Value = 'Education[T.Secondary]'
$ws.Cells.Item(3,2).Value = [double]"-491.0758785794798"
$ws.Cells.Item(3,3).Value = [double]"0.5394962467236093"
$ws.Cells.Item(4,1).Value = 'Education[T.University]'
$ws.Cells.Item(4,2).Value = [double]"539.2362456603779"
$ws.Cells.Item(4,3).Value = [double]"0.4988178776152311"
$ws.Cells.Item(5,1).Value = 'Education[T.Unknown/Other]'
$ws.Cells.Item(5,2).Value = [double]"-2028.321358097021"
$ws.Cells.Item(5,3).Value = [double]"0.09517932393183867"
$ws.Cells.Item(6,1).Value = 'HHSize'
$ws.Cells.Item(6,2).Value = [double]"220.3795222031802"
$ws.Cells.Item(6,3).Value = [double]"0.004254412343829724"
$ws.Cells.Item(7,1).Value = 'Sex'
$ws.Cells.Item(7,2).Value = [double]"-1656.651025088943"
$ws.Cells.Item(7,3).Value = [double]"1.391370274345014e-18"
$ws.Cells.Item(8,1).Value = 'Age'
$ws.Cells.Item(8,2).Value = [double]"-45.15748329531841"
$ws.Cells.Item(8,3).Value = [double]"6.816211247173608e-08"
$ws.Cells.Item(9,1).Value = 'DistSubcenter_res'
$ws.Cells.Item(9,2).Value = [double]"338.9008667333118"
$ws.Cells.Item(9,3).Value = [double]"6.319870694059182e-11"
$ws.Cells.Item(10,1).Value = 'DistCenter_res'
$ws.Cells.Item(10,2).Value = [double]"244.4645755035914"
$ws.Cells.Item(10,3).Value = [double]"3.007743549462916e-11"
$ws.Cells.Item(11,1).Value = 'UrbPopDensity_res'
$ws.Cells.Item(11,2).Value = [double]"0.1577121304786596"
$ws.Cells.Item(11,3).Value = [double]"0.02512127310977947"
$ws.Cells.Item(12,1).Value = 'UrbBuildDensity_res'
$ws.Cells.Item(12,2).Value = [double]"-0.0001434297032114345"
$ws.Cells.Item(12,3).Value = [double]"0.01810547652806192"
$ws.Cells.Item(13,1).Value = 'IntersecDensity_res'
$ws.Cells.Item(13,2).Value = [double]"1.294641344879539"
$ws.Cells.Item(13,3).Value = [double]"0.8534556814127479"
$ws.Cells.Item(14,1).Value = 'street_length_res'
$ws.Cells.Item(14,2).Value = [double]"39.39853574175816"
$ws.Cells.Item(14,3).Value = [double]"6.500153516197236e-11"
$ws.Cells.Item(15,1).Value = 'LU_Comm_res'
$ws.Cells.Item(15,2).Value = [double]"-3092.754430603717"
$ws.Cells.Item(15,3).Value = [double]"0.004990024303608332"
$ws.Cells.Item(16,1).Value = 'LU_UrbFab_res'
$ws.Cells.Item(16,2).Value = [double]"-3236.870871406436"
$ws.Cells.Item(16,3).Value = [double]"2.834099774950492e-06"
$ws.Cells.Item(17,1).Value = 'bike_lane_share_res'
$ws.Cells.Item(17,2).Value = [double]"-2690.921084545842"
$ws.Cells.Item(17,3).Value = [double]"0.1049265661763645"

# --- Sheet 9 ---
$ws = $wb.Worksheets.Item(9)
$ws.Name = 'summ42060716'
$ws.Rows.Item(4).Delete()
$ws.Cells.Item(2,1).Value = 'Intercept'
$ws.Cells.Item(2,2).Value = [double]"9338.145806589659"
$ws.Cells.Item(2,3).Value = [double]"2.690908754874282e-12"
$ws.Cells.Item(3,1).Value = 'Education[T.Secondary]'
$ws.Cells.Item(3,2).Value = [double]"412.5309312811909"
$ws.Cells.Item(3,3).Value = [double]"0.576568427257081"
$ws.Cells.Item(4,1).Value = 'Education[T.University]'
$ws.Cells.Item(4,2).Value = [double]"1365.005319037104"
$ws.Cells.Item(4,3).Value = [double]"0.06350479891333878"
$ws.Cells.Item(5,1).Value = 'Education[T.Unknown/Other]'
$ws.Cells.Item(5,2).Value = [double]"-1711.826162924905"
$ws.Cells.Item(5,3).Value = [double]"0.1408950632315795"
$ws.Cells.Item(6,1).Value = 'HHSize'
$ws.Cells.Item(6,2).Value = [double]"187.9575758739124"
$ws.Cells.Item(6,3).Value = [double]"0.01601337676483489"
$ws.Cells.Item(7,1).Value = 'Sex'
$ws.Cells.Item(7,2).Value = [double]"-1696.172908419456"
$ws.Cells.Item(7,3).Value = [double]"4.330914226337001e-19"
$ws.Cells.Item(8,1).Value = 'Age'
$ws.Cells.Item(8,2).Value = [double]"-37.86348771633716"
$ws.Cells.Item(8,3).Value = [double]"7.340366250535913e-06"
$ws.Cells.Item(9,1).Value = 'DistSubcenter_res'
$ws.Cells.Item(9,2).Value = [double]"389.020306226909"
$ws.Cells.Item(9,3).Value = [double]"1.00004073000788e-13"
$ws.Cells.Item(10,1).Value = 'DistCenter_res'
$ws.Cells.Item(10,2).Value = [double]"243.3761587476874"
$ws.Cells.Item(10,3).Value = [double]"1.045321027768214e-10"
$ws.Cells.Item(11,1).Value = 'UrbPopDensity_res'
$ws.Cells.Item(11,2).Value = [double]"0.1187929155639115"
$ws.Cells.Item(11,3).Value = [double]"0.09645827925620672"
$ws.Cells.Item(12,1).Value = 'UrbBuildDensity_res'
$ws.Cells.Item(12,2).Value = [double]"-0.0001423984064010851"
$ws.Cells.Item(12,3).Value = [double]"0.01831363329217087"
$ws.Cells.Item(13,1).Value = 'IntersecDensity_res'
$ws.Cells.Item(13,2).Value = [double]"-8.558730500614956"
$ws.Cells.Item(13,3).Value = [double]"0.2163359930785129"
$ws.Cells.Item(14,1).Value = 'street_length_res'
$ws.Cells.Item(14,2).Value = [double]"21.37840965606743"
$ws.Cells.Item(14,3).Value = [double]"7.510405107484367e-05"
$ws.Cells.Item(15,1).Value = 'LU_Comm_res'
$ws.Cells.Item(15,2).Value = [double]"-2796.467794362963"
$ws.Cells.Item(15,3).Value = [double]"0.01264920524928372"
$ws.Cells.Item(16,1).Value = 'LU_UrbFab_res'
$ws.Cells.Item(16,2).Value = [double]"-2959.010375824444"
$ws.Cells.Item(16,3).Value = [double]"2.43966386043644e-05"
$ws.Cells.Item(17,1).Value = 'bike_lane_share_res'
$ws.Cells.Item(17,2).Value = [double]"-1819.775414833131"
$ws.Cells.Item(17,3).Value = [double]"0.2812577936904534"
